$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "43.768.52"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "2.296.22"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'108.69"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +11.48%  "
$ws.Range("D6").Value = "'271.64"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").Value = "'0.626"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "'0.617"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("D10").Value = "'47.01"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  +3.66%  "
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("D12").Value = "'8.39"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  +5.13%  "
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").Value = "'15.71"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("D15").Value = "2.636.96"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("D17").Value = "2.291.56"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").Value = "43.790.18"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("D21").Value = "'72.21"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E22").Value = "  +8.86%  "
$ws.Range("D23").Value = "'233.68"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("D24").Value = "'2.95"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  +15.94%  "
$ws.Range("D25").Value = "'9.30"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  -1.83%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'11.34"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").Value = "'40.75"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "  +6.68%  "
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("D31").Value = "'177.91"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("D32").Value = "'21.91"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("D33").Value = "'0.0913"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").Value = "'5.58"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("D35").Value = "'4.88"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  +9.30%  "
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("E37").Value = "  +3.15%  "
$ws.Range("D38").Value = "'0.0359"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("D39").Value = "'3.65"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  +8.52%  "
$ws.Range("D40").Value = "'0.236"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  -3.47%  "
$ws.Range("E41").Value = "  -3.82%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "'66.56"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +5.79%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'1.37"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  -3.50%  "
$ws.Range("D44").Value = "'12.17"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  -1.82%  "
$ws.Range("D45").Value = "'5.49"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  +2.59%  "
$ws.Range("D46").Value = "'8.80"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  -4.41%  "
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("E48").Value = "  +2.19%  "
$ws.Range("D49").Value = "'99.39"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  -1.13%  "
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").Value = "'0.442"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  +5.54%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.53"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  +10.36%  "
